$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values updated to the next day's (automatic) electricity price refresh
$ws.Range("A2").Value = 46013

$ws.Range("B2").Value = 56.16
$ws.Range("C2").Value = 46.5
$ws.Range("D2").Value = 43.72
$ws.Range("E2").Value = 37.63
$ws.Range("F2").Value = 26.09
$ws.Range("G2").Value = 36.4
$ws.Range("H2").Value = 58.66
$ws.Range("I2").Value = 79.36
$ws.Range("J2").Value = 87.26000000000001
$ws.Range("K2").Value = 88.42
$ws.Range("L2").Value = 84.98
$ws.Range("M2").Value = 72.45999999999999
$ws.Range("N2").Value = 66.73999999999999
$ws.Range("O2").Value = 59.92
$ws.Range("P2").Value = 62.03
$ws.Range("Q2").Value = 69.58
$ws.Range("R2").Value = 84.75
$ws.Range("S2").Value = 95.31
$ws.Range("T2").Value = 101.57
$ws.Range("U2").Value = 97.92
$ws.Range("V2").Value = 91.73999999999999
$ws.Range("W2").Value = 85.28
$ws.Range("X2").Value = 84.09999999999999
$ws.Range("Y2").Value = 77.84999999999999
$ws.Range("Z2").Value = 70.59999999999999

$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 94.89

$ws.Range("AD2").Value = 99.73999999999999

$ws.Range("AE2").Value = "16h-18h"
$ws.Range("AF2").Value = 90.03
